$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.754972666666667
$ws.Range("H2").Value = 5.264918
$ws.Range("I2").Value = 0.5110994274238188
$ws.Range("J2").Value = 0.5110994274238188
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 28.72417333333333
$ws.Range("N2").Value = 86.17251999999999
$ws.Range("O2").Value = 0.4233259107972328
$ws.Range("P2").Value = 0.4233259107972328
$ws.Range("Q2").Value = 50.41013907259554
$ws.Range("R2").Value = 453.6912516533599
$ws.Range("S2").Value = 0.2163616306221322
$ws.Range("T2").Value = 0.2163616306221323

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.754972666666667
$ws.Range("H3").Value = 5.264918
$ws.Range("I3").Value = 0.5110994274238188
$ws.Range("J3").Value = 0.5110994274238188
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 30.56986233333333
$ws.Range("N3").Value = 91.709587
$ws.Range("O3").Value = 0.4505269713084062
$ws.Range("P3").Value = 0.4505269713084062
$ws.Range("Q3").Value = 53.64927281876288
$ws.Range("R3").Value = 482.843455368866
$ws.Range("S3").Value = 0.2302640770747137
$ws.Range("T3").Value = 0.2302640770747137

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.754972666666667
$ws.Range("H4").Value = 5.264918
$ws.Range("I4").Value = 0.5110994274238188
$ws.Range("J4").Value = 0.5110994274238188
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.559531999999999
$ws.Range("N4").Value = 25.678596
$ws.Range("O4").Value = 0.126147117894361
$ws.Range("P4").Value = 0.126147117894361
$ws.Range("Q4").Value = 15.02174469945866
$ws.Range("R4").Value = 135.195702295128
$ws.Range("S4").Value = 0.06447371972697286
$ws.Range("T4").Value = 0.06447371972697286

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.678748
$ws.Range("H5").Value = 5.036244
$ws.Range("I5").Value = 0.4889005725761812
$ws.Range("J5").Value = 0.4889005725761812
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 28.72417333333333
$ws.Range("N5").Value = 86.17251999999999
$ws.Range("O5").Value = 0.4233259107972328
$ws.Range("P5").Value = 0.4233259107972328
$ws.Range("Q5").Value = 48.22064853498666
$ws.Range("R5").Value = 433.98583681488
$ws.Range("S5").Value = 0.2069642801751005
$ws.Range("T5").Value = 0.2069642801751005

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.678748
$ws.Range("H6").Value = 5.036244
$ws.Range("I6").Value = 0.4889005725761812
$ws.Range("J6").Value = 0.4889005725761812
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 30.56986233333333
$ws.Range("N6").Value = 91.709587
$ws.Range("O6").Value = 0.4505269713084062
$ws.Range("P6").Value = 0.4505269713084062
$ws.Range("Q6").Value = 51.31909525235866
$ws.Range("R6").Value = 461.871857271228
$ws.Range("S6").Value = 0.2202628942336926
$ws.Range("T6").Value = 0.2202628942336926

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.678748
$ws.Range("H7").Value = 5.036244
$ws.Range("I7").Value = 0.4889005725761812
$ws.Range("J7").Value = 0.4889005725761812
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.559531999999999
$ws.Range("N7").Value = 25.678596
$ws.Range("O7").Value = 0.126147117894361
$ws.Range("P7").Value = 0.126147117894361
$ws.Range("Q7").Value = 14.369297225936
$ws.Range("R7").Value = 129.323675033424
$ws.Range("S7").Value = 0.06167339816738811
$ws.Range("T7").Value = 0.06167339816738812
